$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to remain plain text, matching the inlineStr cell type
    # used in the source workbook (prevents Excel from auto-converting
    # numeric-looking strings like "582.95" into real numbers).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "65.936.47"
$ws.Range("E2").Value = "  -2.29%  "
Set-TextValue $ws.Range("D3") "3.489.34"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "582.95"
$ws.Range("E5").Value = "  -1.38%  "
Set-TextValue $ws.Range("D6") "172.59"
$ws.Range("E6").Value = "  -3.52%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.27%  "
Set-TextValue $ws.Range("D9") "3.488.36"
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("E10").Value = "  -5.48%  "
$ws.Range("E11").Value = "  -1.56%  "
Set-TextValue $ws.Range("D12") "0.411"
$ws.Range("E12").Value = "  -3.50%  "
Set-TextValue $ws.Range("D13") "4.081.61"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("E14").Value = "  +1.24%  "
Set-TextValue $ws.Range("D15") "29.85"
$ws.Range("E15").Value = "  -6.52%  "
Set-TextValue $ws.Range("D16") "65.989.74"
$ws.Range("E16").Value = "  -2.13%  "
Set-TextValue $ws.Range("D17") "0.0000171"
$ws.Range("E17").Value = "  -2.86%  "
Set-TextValue $ws.Range("D18") "3.478.34"
$ws.Range("E18").Value = "  +1.04%  "
Set-TextValue $ws.Range("D19") "5.92"
$ws.Range("E19").Value = "  -3.03%  "
Set-TextValue $ws.Range("D20") "13.91"
$ws.Range("E20").Value = "  -0.42%  "
Set-TextValue $ws.Range("D21") "366.88"
$ws.Range("E21").Value = "  -4.78%  "
Set-TextValue $ws.Range("D22") "7.74"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D23") "72.78"
$ws.Range("E23").Value = "  +2.30%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D24") "1.00"
$ws.Range("E24").Value = "  +0.13%  "
Set-TextValue $ws.Range("D25") "0.0000127"
$ws.Range("E25").Value = "  +6.43%  "
Set-TextValue $ws.Range("D26") "0.534"
$ws.Range("E26").Value = "  +0.55%  "
Set-TextValue $ws.Range("D27") "9.61"
$ws.Range("E27").Value = "  -5.71%  "
Set-TextValue $ws.Range("D28") "0.180"
$ws.Range("E28").Value = "  +2.71%  "
Set-TextValue $ws.Range("D29") "0.999"
$ws.Range("E29").Value = "  -0.07%  "
Set-TextValue $ws.Range("D30") "24.12"
$ws.Range("E30").Value = "  +3.01%  "
Set-TextValue $ws.Range("D31") "5.76"
$ws.Range("E31").Value = "  -4.58%  "
$ws.Range("E32").Value = "  -2.73%  "
Set-TextValue $ws.Range("D33") "0.999"
Set-TextValue $ws.Range("D34") "7.14"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("E35").Value = "  -5.87%  "
Set-TextValue $ws.Range("D36") "1.53"
$ws.Range("E36").Value = "  -1.18%  "
Set-TextValue $ws.Range("D37") "160.51"
$ws.Range("E37").Value = "  -0.43%  "
Set-TextValue $ws.Range("D38") "29.46"
$ws.Range("E38").Value = "  +14.61%  "
Set-TextValue $ws.Range("D39") "0.891"
$ws.Range("E39").Value = "  +1.38%  "
Set-TextValue $ws.Range("D40") "2.827.62"
$ws.Range("E40").Value = "  +4.95%  "
$ws.Range("E41").Value = "  -4.90%  "
Set-TextValue $ws.Range("D42") "6.43"
$ws.Range("E42").Value = "  -2.56%  "
Set-TextValue $ws.Range("D43") "2.56"
$ws.Range("E43").Value = "  -6.27%  "
Set-TextValue $ws.Range("D44") "4.45"
$ws.Range("E44").Value = "  -1.38%  "
Set-TextValue $ws.Range("D45") "0.0682"
$ws.Range("E45").Value = "  -3.93%  "
Set-TextValue $ws.Range("D47") "24.11"
$ws.Range("E47").Value = "  -6.85%  "
$ws.Range("E48").Value = "  -2.75%  "
Set-TextValue $ws.Range("D49") "324.59"
$ws.Range("E49").Value = "  -0.07%  "
Set-TextValue $ws.Range("D50") "0.815"
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("E51").Value = "  -2.56%  "
